# master page framework design
$wb = $excel.ActiveWorkbook

# Rename "Sheet2" to "masterFunction"
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "masterFunction"

# Populate the masterFunction sheet with the master tab labels.
# Enter the sub-tab names first, then the "masterTab" header last, so the
# shared-string table is built in the same order as the authored workbook.
$ws.Cells.Item(2, 1).Value = "Bank"
$ws.Cells.Item(3, 1).Value = "User"
$ws.Cells.Item(4, 1).Value = "Branches"
$ws.Cells.Item(5, 1).Value = "Customers"
$ws.Cells.Item(6, 1).Value = "Brokers"
$ws.Cells.Item(7, 1).Value = "Vendors"
$ws.Cells.Item(8, 1).Value = "Enquiry"
$ws.Cells.Item(9, 1).Value = "Projects"
$ws.Cells.Item(1, 1).Value = "masterTab"

# Fit column A to the new content
$ws.Columns.Item(1).AutoFit()

# Make masterFunction the active/selected sheet
$ws.Activate()
